# Update cryptocurrency price/volume data per latest scrape
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "26.723.85"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +0.15%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.640.48"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -0.20%  "

$ws.Range("E4").Value = "  +0.19%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "217.50"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +1.19%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.504"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -0.09%  "

$ws.Range("E7").Value = "  +0.08%  "

$ws.Range("E8").Value = "  +0.00%  "

$ws.Range("E9").Value = "  -0.04%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "19.15"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +0.48%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.869.75"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -0.10%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "1.651.78"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +1.07%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "4.16"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -0.47%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.527"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -0.40%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "64.58"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -0.73%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "26.718.88"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +0.09%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.0₃0734"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -1.27%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "215.27"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -0.29%  "

$ws.Range("E20").Value = "  +0.06%  "

$ws.Range("E21").Value = "  +0.78%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "2.39"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +6.46%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "6.22"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -0.48%  "

$ws.Range("E24").Value = "  -2.51%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "145.49"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +0.04%  "

$ws.Range("E26").Value = "  -0.03%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.119"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -0.89%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "7.18"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +0.26%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "15.61"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -0.55%  "

$ws.Range("E30").Value = "  -1.08%  "

$ws.Range("E31").Value = "  +1.16%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.38"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +0.61%  "

$ws.Range("E33").Value = "  -0.54%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.284.53"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +0.30%  "

$ws.Range("E35").Value = "  -0.07%  "

$ws.Range("E36").Value = "  +1.11%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.0176"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -1.00%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.538"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +1.33%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.817"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -1.37%  "

$ws.Range("E40").Value = "  +0.05%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.805"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -0.95%  "

$ws.Range("E42").Value = "  -1.79%  "

$ws.Range("E43").Value = "  -2.65%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "1.779.75"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -0.12%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "61.23"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +3.53%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "91.86"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +0.32%  "

$ws.Range("E47").Value = "  +0.47%  "

$ws.Range("E48").Value = "  +0.24%  "

$ws.Range("E49").Value = "  -1.41%  "

$ws.Range("E50").Value = "  +0.25%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.407"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -0.04%  "

